# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "66.919.97"
Set-TextValue $ws.Range("E2") "  +2.07%  "

Set-TextValue $ws.Range("D3") "3.436.61"
Set-TextValue $ws.Range("E3") "  +1.22%  "

Set-TextValue $ws.Range("E4") "  -0.11%  "

Set-TextValue $ws.Range("D5") "575.73"
Set-TextValue $ws.Range("E5") "  +2.75%  "

Set-TextValue $ws.Range("D6") "186.79"
Set-TextValue $ws.Range("E6") "  +6.13%  "

Set-TextValue $ws.Range("E7") "  -0.08%  "

Set-TextValue $ws.Range("D8") "3.430.05"
Set-TextValue $ws.Range("E8") "  +1.32%  "

Set-TextValue $ws.Range("E9") "  -0.08%  "

Set-TextValue $ws.Range("E10") "  -0.10%  "

Set-TextValue $ws.Range("D11") "0.639"
Set-TextValue $ws.Range("E11") "  +0.34%  "

Set-TextValue $ws.Range("D12") "57.52"
Set-TextValue $ws.Range("E12") "  +7.03%  "

Set-TextValue $ws.Range("E13") "  -1.17%  "

Set-TextValue $ws.Range("E14") "  +2.08%  "

Set-TextValue $ws.Range("D15") "3.993.29"
Set-TextValue $ws.Range("E15") "  +1.32%  "

Set-TextValue $ws.Range("D16") "18.86"
Set-TextValue $ws.Range("E16") "  +2.97%  "

Set-TextValue $ws.Range("D17") "3.440.60"
Set-TextValue $ws.Range("E17") "  +1.72%  "

Set-TextValue $ws.Range("D18") "66.954.15"
Set-TextValue $ws.Range("E18") "  +2.27%  "

Set-TextValue $ws.Range("E19") "  -0.55%  "

Set-TextValue $ws.Range("D20") "11.99"
Set-TextValue $ws.Range("E20") "  +1.05%  "

Set-TextValue $ws.Range("E21") "  +1.23%  "

Set-TextValue $ws.Range("D22") "487.83"
Set-TextValue $ws.Range("E22") "  +5.11%  "

Set-TextValue $ws.Range("E23") "  +12.18%  "

Set-TextValue $ws.Range("D24") "16.87"
Set-TextValue $ws.Range("E24") "  +17.55%  "

Set-TextValue $ws.Range("D25") "4.30"
Set-TextValue $ws.Range("E25") "  +3.85%  "

Set-TextValue $ws.Range("D26") "89.26"
Set-TextValue $ws.Range("E26") "  +2.24%  "

Set-TextValue $ws.Range("E27") "  +0.86%  "

Set-TextValue $ws.Range("D28") "10.90"
Set-TextValue $ws.Range("E28") "  +1.68%  "

Set-TextValue $ws.Range("D29") "8.96"
Set-TextValue $ws.Range("E29") "  +2.44%  "

Set-TextValue $ws.Range("D30") "31.07"
Set-TextValue $ws.Range("E30") "  -0.01%  "

Set-TextValue $ws.Range("E31") "  +12.17%  "

Set-TextValue $ws.Range("D32") "603.01"
Set-TextValue $ws.Range("E32") "  +4.21%  "

Set-TextValue $ws.Range("D33") "64.67"
Set-TextValue $ws.Range("E33") "  +1.96%  "

Set-TextValue $ws.Range("D34") "11.76"
Set-TextValue $ws.Range("E34") "  +2.21%  "

Set-TextValue $ws.Range("E36") "  -0.10%  "

Set-TextValue $ws.Range("E37") "  +2.00%  "

Set-TextValue $ws.Range("B38") "InjectiveProtocol"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D38") "36.60"
Set-TextValue $ws.Range("E38") "  +1.78%  "

Set-TextValue $ws.Range("B39") "PEPE"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D39") "0.0₃0777"
Set-TextValue $ws.Range("E39") "  +4.78%  "

Set-TextValue $ws.Range("D40") "0.384"
Set-TextValue $ws.Range("E40") "  +2.43%  "

Set-TextValue $ws.Range("D41") "3.43"
Set-TextValue $ws.Range("E41") "  -4.66%  "

Set-TextValue $ws.Range("D42") "3.181.32"
Set-TextValue $ws.Range("E42") "  +1.85%  "

Set-TextValue $ws.Range("E43") "  +2.15%  "

Set-TextValue $ws.Range("D44") "0.0427"
Set-TextValue $ws.Range("E44") "  +2.19%  "

Set-TextValue $ws.Range("E45") "  +4.38%  "

Set-TextValue $ws.Range("E46") "  +1.16%  "

Set-TextValue $ws.Range("E47") "  +0.76%  "

Set-TextValue $ws.Range("E48") "  +14.21%  "

Set-TextValue $ws.Range("E49") "  -0.11%  "

Set-TextValue $ws.Range("B50") "Monero"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D50") "141.28"
Set-TextValue $ws.Range("E50") "  +0.52%  "

Set-TextValue $ws.Range("B51") "THORChain"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "8.60"
Set-TextValue $ws.Range("E51") "  +1.60%  "

